$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.297.98"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.559.22"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'606.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'144.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "3.558.31"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "4.163.10"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "'30.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "3.551.86"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "66.342.69"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "'11.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("D20").Value = "'6.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'14.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'431.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "'79.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").Value = "3.700.14"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'0.0000119"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'9.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").Value = "'7.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("B32").Value = "RenzoRestakedETH"
$ws.Range("C32").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D32").Value = "3.553.81"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'25.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("E35").Value = "  -5.68%  "
$ws.Range("D36").Value = "'7.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'175.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("D41").Value = "'0.0849"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").Value = "'0.888"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").Value = "'1.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").Value = "'46.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  +5.16%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "'25.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").Value = "'7.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "'23.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.94%  "
